$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the course number string in F3 (00000002 -> 00000004), keeping it as text
$ws.Range("F3").Value = "'00000004"

# Remove the two student rows (0000000001/zhazha/88.5 and 0000000002/hahaha/77.5).
# Deleting row 6 twice shifts the remaining rows (old 8, 9) up into rows 6, 7.
$ws.Rows(6).Delete()
$ws.Rows(6).Delete()

# The surviving rows (now 6 and 7) keep their student id / name but lose their score.
$ws.Range("C6:C7").ClearContents()
